$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format risky numeric-looking Price cells as Text so the
# trailing/insignificant digits and string type are preserved,
# matching the source data (these are display strings, not numbers).
$textCells = @("D5", "D6", "D12", "D14", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D31", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D44", "D47", "D51")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.195.63'
$ws.Range("E2").Value = '  +4.71%  '
$ws.Range("D3").Value = '2.365.90'
$ws.Range("E3").Value = '  +3.56%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '520.01'
$ws.Range("E5").Value = '  +3.34%  '
$ws.Range("D6").Value = '135.10'
$ws.Range("E6").Value = '  +3.92%  '
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  +2.03%  '
$ws.Range("D9").Value = '2.363.08'
$ws.Range("E9").Value = '  +2.96%  '
$ws.Range("E10").Value = '  +8.24%  '
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = '5.24'
$ws.Range("E12").Value = '  +6.56%  '
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").Value = '23.92'
$ws.Range("E14").Value = '  +3.16%  '
$ws.Range("D15").Value = '2.786.11'
$ws.Range("E15").Value = '  +3.47%  '
$ws.Range("D16").Value = '57.084.60'
$ws.Range("E16").Value = '  +4.49%  '
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").Value = '  +3.75%  '
$ws.Range("D18").Value = '2.360.94'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").Value = '10.57'
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("D20").Value = '4.29'
$ws.Range("E20").Value = '  +3.12%  '
$ws.Range("D21").Value = '323.54'
$ws.Range("E21").Value = '  +5.60%  '
$ws.Range("D22").Value = '6.74'
$ws.Range("E22").Value = '  +6.00%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '61.50'
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("D25").Value = '0.996'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("E26").Value = '  +6.56%  '
$ws.Range("E27").Value = '  +5.33%  '
$ws.Range("D28").Value = '171.40'
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").Value = '0.0₃0744'
$ws.Range("E29").Value = '  +5.06%  '
$ws.Range("E30").Value = '  +9.71%  '
$ws.Range("D31").Value = '6.29'
$ws.Range("E31").Value = '  +4.12%  '
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("E33").Value = '  +2.42%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.39%  '
$ws.Range("D36").Value = '0.957'
$ws.Range("E36").Value = '  +2.02%  '
$ws.Range("D37").Value = '1.26'
$ws.Range("E37").Value = '  +5.05%  '
$ws.Range("D38").Value = '4.03'
$ws.Range("E38").Value = '  +7.24%  '
$ws.Range("E39").Value = '  +7.40%  '
$ws.Range("D40").Value = '37.57'
$ws.Range("E40").Value = '  +3.72%  '
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("D42").Value = '140.66'
$ws.Range("E42").Value = '  +12.06%  '
$ws.Range("E43").Value = '  +5.41%  '
$ws.Range("D44").Value = '278.74'
$ws.Range("E44").Value = '  +12.89%  '
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("E46").Value = '  +3.69%  '
$ws.Range("D47").Value = '0.0930'
$ws.Range("E47").Value = '  +3.46%  '
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("E50").Value = '  +4.41%  '
$ws.Range("D51").Value = '17.03'
$ws.Range("E51").Value = '  +3.30%  '

# Clear the temporary number-format override so the cell style
# index returns to the sheet default (unstyled), matching source.
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).ClearFormats()
}
